$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3320
$ws.Range("I76").Value = 3304.348
$ws.Range("J76").Value = 3500
$ws.Range("K76").Value = 3304.348
$ws.Range("L76").Value = 3500
$ws.Range("M76").Value = -2989.348
$ws.Range("N76").Value = -4130

$ws.Range("H79").Value = 3320
$ws.Range("I79").Value = 3304.348
$ws.Range("J79").Value = 3500
$ws.Range("K79").Value = 3304.348
$ws.Range("L79").Value = 3500
$ws.Range("M79").Value = -2212.348
$ws.Range("N79").Value = -5684

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 4272.8335
$ws.Range("I74").Value = 1086.3077
$ws.Range("J74").Value = 6073.913
$ws.Range("K74").Value = 1086.3077
$ws.Range("L74").Value = 6073.913
$ws.Range("M74").Value = -212.3077000000001
$ws.Range("N74").Value = -7821.913

$ws.Range("H77").Value = 4272.8335
$ws.Range("I77").Value = 1086.3077
$ws.Range("J77").Value = 6073.913
$ws.Range("K77").Value = 5431.538500000001
$ws.Range("L77").Value = 30369.565
$ws.Range("M77").Value = -1063.538500000001
$ws.Range("N77").Value = -39105.565

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1732.8462
$ws.Range("I99").Value = 1728.5
$ws.Range("J99").Value = 1739.8
$ws.Range("K99").Value = 1728.5
$ws.Range("L99").Value = 1739.8
$ws.Range("M99").Value = -230.5
$ws.Range("N99").Value = -4735.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 12999.667
$ws.Range("J28").Value = 12999.667
$ws.Range("L28").Value = 12999.667
$ws.Range("N28").Value = -13489.667

$ws.Range("H31").Value = 1637.3334
$ws.Range("I31").Value = 1205.619
$ws.Range("J31").Value = 2392.8333
$ws.Range("K31").Value = 1205.619
$ws.Range("L31").Value = 2392.8333
$ws.Range("M31").Value = -910.6189999999999
$ws.Range("N31").Value = -2982.8333

$ws.Range("H34").Value = 1637.3334
$ws.Range("I34").Value = 1205.619
$ws.Range("J34").Value = 2392.8333
$ws.Range("K34").Value = 1205.619
$ws.Range("L34").Value = 2392.8333
$ws.Range("M34").Value = -1003.619
$ws.Range("N34").Value = -2796.8333

$ws.Range("H35").Value = 1209.3
$ws.Range("I35").Value = 1209.3
$ws.Range("K35").Value = 1209.3
$ws.Range("M35").Value = -915.3

$ws.Range("H58").Value = 3459.2104
$ws.Range("I58").Value = 2391.3635
$ws.Range("J58").Value = 3894.2593
$ws.Range("K58").Value = 2391.3635
$ws.Range("L58").Value = 3894.2593
$ws.Range("M58").Value = -2188.3635
$ws.Range("N58").Value = -4300.2593

$ws.Range("H136").Value = 3459.2104
$ws.Range("I136").Value = 2391.3635
$ws.Range("J136").Value = 3894.2593
$ws.Range("K136").Value = 7174.0905
$ws.Range("L136").Value = 11682.7779
$ws.Range("M136").Value = -4624.0905
$ws.Range("N136").Value = -16782.7779

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()

$ws.Range("H34").Value = 28573048
$ws.Range("I34").Value = 199
$ws.Range("J34").Value = 33335190
$ws.Range("K34").Value = 597
$ws.Range("L34").Value = 100005570
$ws.Range("M34").Value = -513
$ws.Range("N34").Value = -100005738

$ws.Range("H39").Value = 3318.182
$ws.Range("I39").Value = 4000
$ws.Range("J39").Value = 3250
$ws.Range("K39").Value = 12000
$ws.Range("L39").Value = 9750
$ws.Range("M39").Value = -11706
$ws.Range("N39").Value = -10338

$ws.Range("H52").Value = 681.1111
$ws.Range("J52").Value = 681.1111
$ws.Range("L52").Value = 2043.3333
$ws.Range("N52").Value = -2575.3333

$ws.Range("H55").Value = 208.80952
$ws.Range("I55").Value = 105
$ws.Range("J55").Value = 260.7143
$ws.Range("K55").Value = 315
$ws.Range("L55").Value = 782.1428999999999
$ws.Range("M55").Value = -138
$ws.Range("N55").Value = -1136.1429

$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()

$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()

$ws.Range("H131").Value = 1756921.6
$ws.Range("I131").Value = 5885.5557
$ws.Range("J131").Value = 2085240.9
$ws.Range("K131").Value = 17656.6671
$ws.Range("L131").Value = 6255722.699999999
$ws.Range("M131").Value = -12616.6671
$ws.Range("N131").Value = -6265802.699999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 136665.67
$ws.Range("I80").Value = 2500
$ws.Range("J80").Value = 174998.72
$ws.Range("K80").Value = 2500
$ws.Range("L80").Value = 174998.72
$ws.Range("M80").Value = -1502
$ws.Range("N80").Value = -176994.72

$ws.Range("H83").Value = 136665.67
$ws.Range("I83").Value = 2500
$ws.Range("J83").Value = 174998.72
$ws.Range("K83").Value = 12500
$ws.Range("L83").Value = 874993.6
$ws.Range("M83").Value = -7508
$ws.Range("N83").Value = -884977.6

$ws.Range("H107").Value = 567.9524
$ws.Range("I107").Value = 244
$ws.Range("J107").Value = 697.5333000000001
$ws.Range("K107").Value = 244
$ws.Range("L107").Value = 697.5333000000001
$ws.Range("M107").Value = 1676
$ws.Range("N107").Value = -4537.5333

$ws.Range("H132").Value = 2761.4443
$ws.Range("I132").Value = 1955.2106
$ws.Range("J132").Value = 4676.25
$ws.Range("K132").Value = 5865.6318
$ws.Range("L132").Value = 14028.75
$ws.Range("M132").Value = -3335.6318
$ws.Range("N132").Value = -19088.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2636.8
$ws.Range("I16").Value = 2760.8572
$ws.Range("J16").Value = 900
$ws.Range("K16").Value = 2760.8572
$ws.Range("L16").Value = 900
$ws.Range("M16").Value = -2590.8572
$ws.Range("N16").Value = -1240

$ws.Range("H22").Value = 368
$ws.Range("I22").Value = 313.33334
$ws.Range("J22").Value = 450
$ws.Range("K22").Value = 313.33334
$ws.Range("L22").Value = 450
$ws.Range("M22").Value = -18.33334000000002
$ws.Range("N22").Value = -1040

$ws.Range("H27").Value = 368
$ws.Range("I27").Value = 313.33334
$ws.Range("J27").Value = 450
$ws.Range("K27").Value = 313.33334
$ws.Range("L27").Value = 450
$ws.Range("M27").Value = -206.33334
$ws.Range("N27").Value = -664

$ws.Range("H132").Value = 6950.32
$ws.Range("I132").Value = 8818.75
$ws.Range("J132").Value = 3628.6667
$ws.Range("K132").Value = 26456.25
$ws.Range("L132").Value = 10886.0001
$ws.Range("M132").Value = -23926.25
$ws.Range("N132").Value = -15946.0001
